$wb = $excel.ActiveWorkbook

# --- Main sheet data updates ---
$wsMain = $wb.Worksheets.Item("Main")

# Price
$wsMain.Range("M2").Value = 400
# Shares
$wsMain.Range("M3").Value = 64
# Quarter label next to Shares (Q124 -> Q224)
$wsMain.Range("N3").Value = "Q224"

# MC (M4) keeps its existing formula +M2*M3 and recalculates automatically

# Cash
$wsMain.Range("M5").Value = 1670
$wsMain.Range("N5").Value = "Q224"
$wsMain.Range("N5").HorizontalAlignment = -4152

# Debt
$wsMain.Range("M6").Value = 2174
$wsMain.Range("N6").Value = "Q224"
$wsMain.Range("N6").HorizontalAlignment = -4152

# EV
$wsMain.Range("M7").Formula = "=+M4-M5+M6"

# --- View / selection state ---
$wsModel = $wb.Worksheets.Item("Model")
$wsModel.Activate()
$wsModel.Range("B34").Select()

$wsMain.Activate()
$wsMain.Range("B3").Select()

Write-Output "done"
